$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the (old) second column. "StageID" stays in
# column A; the new column B will hold "IsOpen" and a further new column C
# will hold "MaxStart". Everything that used to live in columns B.. shifts
# one column to the right.
$ws.Columns("B").Insert()
$ws.Columns("B").ColumnWidth = 25

# --- formatting: clone column A's per-row styles onto the two new columns
$xlPasteFormats = -4122
$ws.Range("A1:A25").Copy()
$ws.Range("B1:B25").PasteSpecial($xlPasteFormats)
$ws.Range("A1:A25").Copy()
$ws.Range("C1:C25").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# --- row 1 header ("All" spans all three columns now)
$ws.Range("B1").Value = "All"
$ws.Range("C1").Value = "All"

# --- row 2 column headers
$ws.Range("B2").Value = "IsOpen"
$ws.Range("C2").Value = "MaxStart"

# --- data rows 3-25: B = IsOpen (only StageID 0, i.e. row 3, is open),
#     C = MaxStart (constant 3 for every stage)
for ($r = 3; $r -le 25; $r++) {
    if ($r -eq 3) {
        $ws.Cells.Item($r, 2).Value = 1
    } else {
        $ws.Cells.Item($r, 2).Value = 0
    }
    $ws.Cells.Item($r, 3).Value = 3
}

$ws.Range("B2").Select() | Out-Null
